$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every Price cell we touch so values stay text
# (matching the original inlineStr cells) instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.432.90"
$ws.Range("D3").Value = "1.848.96"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "240.72"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "0.6266"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D8").Value = "0.07684"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "0.2912"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").Value = "24.78"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").Value = "0.07748"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.854.27"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "5.025"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "0.6813"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "0.00001074"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "83.54"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "6.167"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "29.455.33"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "228.38"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "7.406"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "157.25"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").Value = "8.396"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "17.69"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "1.343"
$ws.Range("E28").Value = "  +5.02%  "
$ws.Range("D29").Value = "1.466"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").Value = "0.05644"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").Value = "4.116"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "4.032"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "1.841"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "0.7024"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "2.772"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.228.56"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").Value = "0.01787"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").Value = "6.533"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "0.9039"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D43").Value = "1.992.01"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "101.70"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "65.86"
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.150"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000119"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "0.4010"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("D50").Value = "9.000"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  +0.13%  "
